$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: cardholder first/last name
$ws.Range("C2").Value = "Hartmut"
$ws.Range("C3").Value = "Mohaupt"

# Card number (16-digit numeric-looking string must stay TEXT, not become a
# number - otherwise Excel would render it in scientific notation). Enter it
# as a text-producing formula (forces the string type) then collapse the
# formula down to a plain static value via a values-only self paste, which
# also keeps the original "s=8" cell style untouched.
$ws.Range("B3").Formula = "=""2570314725427075"""
$ws.Range("B3").Copy()
$ws.Range("B3").PasteSpecial(-4163)

# Opening balance banner date
$ws.Range("D5").Value = "KONTOSTAND AM 26.12.2024"

# Row 6 transaction
$ws.Range("B6").Value = "27.12."
$ws.Range("C6").Value = "28.12."
$ws.Range("D6").Value = "KARTENZ./27.12 EDEKA RO"
$ws.Range("E6").Value = "79,35-"

# Row 7 transaction
$ws.Range("B7").Value = "28.12."
$ws.Range("C7").Value = "29.12."
$ws.Range("D7").Value = "BEITRAG Allianz SE K-67236422"
$ws.Range("E7").Value = "52,91-"

# Row 8 transaction
$ws.Range("B8").Value = "29.12."
$ws.Range("C8").Value = "30.12."
$ws.Range("D8").Value = "AMAZON.DE MKTPLC EU JCCDXP"
$ws.Range("E8").Value = "86,83-"

# Rows 9 and 10 - the last two transactions were dropped from the
# statement entirely, leaving blank rows (matching the blank spacer
# row 11's look). Pull row 11's cell formatting onto E9/E10 first
# (cheap single-step style reuse), then nudge E9 back to the
# center-aligned variant, and finally wipe all the text out.
$ws.Range("E11").Copy()
$ws.Range("E9").PasteSpecial(-4122)
$ws.Range("E9").HorizontalAlignment = -4108
$ws.Range("E10").PasteSpecial(-4122)

$ws.Range("B9:E9").Value = ""
$ws.Range("B10:E10").Value = ""

# Closing balance banner + amount
$ws.Range("D12").Value = "KONTOSTAND AM 02.01.2025"
$ws.Range("E12").Value = "219,09-"

# Next statement date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 07.01.2025"
